$wb = $excel.ActiveWorkbook

# --- testdata sheet: update the stored Stripe customer id (A15) ---
$wsData = $wb.Worksheets.Item("testdata")
$wsData.Range("A15").Value = "cus_OdVhGd9iTYfzfB"

# --- cxCreationValidKey sheet: widen column B to fit its (now longer) hyperlink text ---
$wsValid = $wb.Worksheets.Item("cxCreationValidKey")
$wsValid.Columns.Item(2).ColumnWidth = 18.65

# --- make cxCreationValidKey the active/selected sheet (was "testdata") ---
$wsValid.Activate()
